$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.279.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.33%  '
$ws.Range("D3").Value = "'1.610.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.75%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = "'213.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.53%  '
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("E7").Value = '  +0.42%  '
$ws.Range("D8").Value = "'0.250"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.10%  '
$ws.Range("E9").Value = '  +0.43%  '
$ws.Range("D10").Value = "'18.45"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.36%  '
$ws.Range("E11").Value = '  +0.49%  '
$ws.Range("D12").Value = "'1.834.60"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.70%  '
$ws.Range("D13").Value = "'1.604.27"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.37%  '
$ws.Range("E14").Value = '  +0.18%  '
$ws.Range("E15").Value = '  +1.30%  '
$ws.Range("D16").Value = "'26.263.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.24%  '
$ws.Range("D17").Value = "'62.10"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.21%  '
$ws.Range("E18").Value = '  +0.86%  '
$ws.Range("E19").Value = '  -0.11%  '
$ws.Range("D20").Value = "'201.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.39%  '
$ws.Range("E21").Value = '  +1.22%  '
$ws.Range("E22").Value = '  +0.73%  '
$ws.Range("E23").Value = '  +0.93%  '
$ws.Range("E24").Value = '  +4.28%  '
$ws.Range("D25").Value = "'143.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.47%  '
$ws.Range("E26").Value = '  -0.08%  '
$ws.Range("D27").Value = "'0.122"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.66%  '
$ws.Range("D28").Value = "'15.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.78%  '
$ws.Range("E29").Value = '  +2.37%  '
$ws.Range("D30").Value = "'0.0500"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.76%  '
$ws.Range("E31").Value = '  +0.35%  '
$ws.Range("E32").Value = '  +2.70%  '
$ws.Range("E33").Value = '  -0.21%  '
$ws.Range("E34").Value = '  +1.33%  '
$ws.Range("E35").Value = '  +0.57%  '
$ws.Range("D36").Value = "'1.156.23"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.80%  '
$ws.Range("E37").Value = '  +0.95%  '
$ws.Range("E38").Value = '  -0.11%  '
$ws.Range("E39").Value = '  +1.24%  '
$ws.Range("E40").Value = '  +0.17%  '
$ws.Range("D41").Value = "'0.496"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.37%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = "'5.34"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.11%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = "'0.784"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.40%  '
$ws.Range("D44").Value = "'1.745.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.62%  '
$ws.Range("D45").Value = "'92.71"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.17%  '
$ws.Range("D47").Value = "'1.53"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.28%  '
$ws.Range("D48").Value = "'53.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.91%  '
$ws.Range("D49").Value = "'0.0508"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.80%  '
$ws.Range("E50").Value = '  -0.25%  '
$ws.Range("E51").Value = '  -0.28%  '
